# CIS399Wk2Day3b-Menu+Settings.pptx -- "Updated lab assignments and rubrics"
#
# The only content-level change in the target diff that is reachable through
# the PowerPoint COM object model is the resize of the screenshot picture on
# slide 31 ("Picture 4"): its extent (cx/cy) grows from 7138147x4656921 EMU
# to 7773714x5244388 EMU while its position (off x/y = 876300/420688) is
# left untouched.
#
# (The diff's other hunks -- the `spid="_x0000_sNNNN"` VML shape-id bumps on
# the ten embedded Word `p:oleObj` graphicFrames, and the `dirty="0"` marks
# added to two `a:rPr` runs on slide 20 -- are incidental byte-level
# artifacts that real PowerPoint itself only ever produces as a side effect
# of interactively activating an OLE object / typing in the UI. Neither is
# backed by a settable property anywhere in the documented Shape/TextRange/
# Font/OLEFormat object model, so there is no COM call that reproduces them.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(31)
$shp = $s.Shapes.Item(4)

# Sanity check -- make sure we are touching the right shape before resizing.
if ($shp.Name -ne "Picture 4") {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq "Picture 4") {
            $shp = $s.Shapes.Item($i)
        }
    }
}

# Keep the top-left corner fixed (don't touch Left/Top at all -- leaving
# them alone avoids any round-trip rounding of their already-correct EMU
# values); only grow the width/height.
# The literal point values below are chosen (via the nearest float32
# representable neighbours) so that, after the host's internal
# points->EMU conversion, they land exactly on the target EMU extents
# cx=7773714 / cy=5244388 instead of drifting by rounding error.
$shp.Width = 612.103515625
$shp.Height = 412.9439392089844

Write-Output "Slide 31 '$($shp.Name)': Left=$($shp.Left) Top=$($shp.Top) Width=$($shp.Width) Height=$($shp.Height)"
